# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   for every row that shares that status string (Overview + per-locale sheets).
# - Each locale sheet's "Latest Handback DateTime" is refreshed to the handback run time.
# - The stale "handback file is not the latest" Error Detail messages are cleared
#   now that the handback is in sync.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

$newStatus = "Handed back: in sync with en-US"

# --- zh-cn sheet (table columns: C=Status, K=Latest Handback DateTime, P=Error Detail) ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-10-26 07:29:42"
$wsZhCn.Range("K3").Value = "2016-10-26 07:29:42"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet (table columns: C=Status, K=Latest Handback DateTime, P=Error Detail) ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-10-26 07:29:59"
$wsDeDe.Range("K3").Value = "2016-10-26 07:29:59"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

# --- Overview sheet mirrors the status for both locales ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- Column widths widen to fit the longer status text / shrink now Error Detail is empty ---
$wsOverview.Range("E:E").ColumnWidth = 29.9777050018311
$wsOverview.Range("F:F").ColumnWidth = 29.9777050018311

$wsZhCn.Range("C:C").ColumnWidth = 29.9777050018311
$wsZhCn.Range("P:P").ColumnWidth = 13.7470531463623

$wsDeDe.Range("C:C").ColumnWidth = 29.9777050018311
$wsDeDe.Range("P:P").ColumnWidth = 13.7470531463623
